$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (2007 data) - shifts remaining rows (2010, 2012, 2015, 2017) up by one
$ws.Rows.Item(2).Delete()
